$d = $word.ActiveDocument

# 1) Update the letter date.
$d.Content.Find.Execute(
    "September 19, 2025", $true, $false, $false, $false, $false,
    $true, 1, $false, "September 21, 2025", 2) | Out-Null

# 2) Split the single-line mailing address into two paragraphs:
#    "1730 Highland Place" and a new paragraph "Berkeley, CA 94709".
$d.Content.Find.Execute(
    "1730 Highland Place, Berkeley CA 94709", $true, $false, $false, $false, $false,
    $true, 1, $false, "1730 Highland Place^pBerkeley, CA 94709", 2) | Out-Null

# Give the freshly split-off "Berkeley, CA 94709" paragraph the same
# run formatting (Arial, 11pt / sz 22) as the rest of the address block.
$addrPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "Berkeley, CA 94709`r") {
        $addrPara = $p
        break
    }
}
if ($addrPara -ne $null) {
    $addrPara.Range.Font.Name = "Arial"
    $addrPara.Range.Font.NameAscii = "Arial"
    $addrPara.Range.Font.NameBi = "Arial"
    $addrPara.Range.Font.Size = 11
    $addrPara.Range.Font.SizeBi = 11
}

# 3) Remove the now-redundant blank "No Spacing" paragraph that used to
#    follow "Board of Directors".
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Board of Directors") {
        $next = $d.Paragraphs.Item($i + 1)
        if ($next.Range.Text -eq "`r" -and $next.Style.NameLocal -eq "No Spacing") {
            $next.Range.Delete()
        }
        break
    }
}

Write-Output "Done: edited date, split mailing address, removed blank paragraph."
